# Trade #78 closed at 2026-02-16 21:36:27 - leadlag DOWN +0.000%
#
# This script:
#   1. Updates the Summary sheet roll-up numbers.
#   2. Closes trade #53 (leadlag sheet row 43 / "All Trades" sheet row 54)
#      by filling in its exit price / status / P&L / exit reason / duration,
#      and mirrors the now-closed trade into the "All Trades" sheet as a new
#      row.
#   3. Appends a brand-new OPEN trade (#78) to the "leadlag" sheet (row 58).
#   4. Updates the Comparison sheet's leadlag summary stats.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a string value into a cell WITHOUT letting Excel's
# automatic type-sniffing turn it into a number/date/percentage. We do this
# by staging the literal text into a scratch cell that's explicitly
# formatted as Text ("@"), then copying just the VALUE (PasteSpecial
# xlPasteValues = -4163) into the destination, which keeps the destination
# cell's style untouched (General) while preserving the literal text.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($Sheet, [string]$Address, [string]$Text)

    $scratch = $Sheet.Range("ZZ999")
    $scratch.NumberFormat = "@"
    $scratch.Value = $Text
    $scratch.Copy()
    $Sheet.Range($Address).PasteSpecial(-4163)
    $scratch.Clear()
}

# ---------------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("C2").Value = 53
Set-TextValue $wsSummary "D2" "66.0%"
Set-TextValue $wsSummary "E2" "+13.3869%"
Set-TextValue $wsSummary "F2" "+0.2526%"
Set-TextValue $wsSummary "E3" "+9.1853%"
Set-TextValue $wsSummary "F3" "+0.1640%"

# ---------------------------------------------------------------------------
# 2. leadlag sheet - close trade #53 (row 43)
# ---------------------------------------------------------------------------
$wsLeadlag = $wb.Worksheets.Item("leadlag")

$wsLeadlag.Range("G43").Value = 68875.76079499999
Set-TextValue $wsLeadlag "H43" "CLOSED"
$wsLeadlag.Range("I43").Value = -0.2518
$wsLeadlag.Range("J43").Value = -2.52
Set-TextValue $wsLeadlag "M43" "time_exit_5min"
$wsLeadlag.Range("N43").Value = 5

# Mirror the now-closed trade #53 into the "All Trades" sheet as new row 54.
# ("All Trades" numbers rows sequentially by closed-trade count, which here
# happens to also be 53 - the 53rd closed trade overall.)
$wsAllTrades = $wb.Worksheets.Item("All Trades")
$wsLeadlag.Range("A43:N43").Copy($wsAllTrades.Range("A54:N54"))
$wsAllTrades.Range("A54").Value = 53

# ---------------------------------------------------------------------------
# 3. leadlag sheet - append new OPEN trade #78 (row 58)
# ---------------------------------------------------------------------------
# Seed row 58 from row 57 (same shape: leadlag / OPEN trade) via Copy so the
# date/time text cells stay plain text instead of being auto-converted to
# Excel date/time serials, then overwrite the cells that actually differ.
$wsLeadlag.Range("A57:N57").Copy($wsLeadlag.Range("A58:N58"))

$wsLeadlag.Range("A58").Value = 78
Set-TextValue $wsLeadlag "C58" "21:36:27"
Set-TextValue $wsLeadlag "E58" "DOWN"
$wsLeadlag.Range("F58").Value = 68615.23
Set-TextValue $wsLeadlag "L58" "Binance leading with -0.148% move"

# ---------------------------------------------------------------------------
# 4. Comparison sheet
# ---------------------------------------------------------------------------
$wsComparison = $wb.Worksheets.Item("Comparison")

Set-TextValue $wsComparison "D2" "2.76"
Set-TextValue $wsComparison "F2" "-0.3266%"
Set-TextValue $wsComparison "G2" "1.70"
